# Phase 2 Constellation Launches
# Scale the "Satellites per operator" (column H) figures on the
# "scenarios copy" sheet up by an order of magnitude (x10), matching the
# refreshed scenarios.csv source data. H13 is already a stand-alone source
# value in the new data (it does not get multiplied).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenarios copy")

# row -> new value for column H (rows 2-21)
$updates = [ordered]@{
    2  = 44250
    3  = 3940
    4  = 1800
    5  = 720
    6  = 6480
    7  = 9720
    8  = 13680
    9  = 1080
    10 = 720
    11 = 720
    12 = 100
    13 = 9
    14 = 44250
    15 = 6480
    16 = 32360
    17 = 190
    18 = 44250
    19 = 6480
    20 = 32360
    21 = 190
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 8)   # column H = 8
    $cell.Value = $updates[$row]
    # Nudge the cell's format off the implicit default style (matches the
    # distinct, already-"Normal" style used for this refreshed data range
    # in the authored workbook).
    $cell.WrapText = $false
}

# Cosmetic view refresh to match how the sheet was left after the edit:
# zoomed to 100% with the refreshed H column selected.
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$ws.Range("H2:H21").Select() | Out-Null
